$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Update header labels on existing sheets
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# 2. Add the new "PO Forecast" worksheet after "Monthly Trend"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# match the page margins used by the other sheets in the workbook
$ws3.PageSetup.LeftMargin = $ws1.PageSetup.LeftMargin
$ws3.PageSetup.RightMargin = $ws1.PageSetup.RightMargin
$ws3.PageSetup.TopMargin = $ws1.PageSetup.TopMargin
$ws3.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$ws3.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$ws3.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

# 3. Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# copy the bold/border/centered header style from an existing header cell
$ws1.Range("B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# 4. Forecast data rows
$ws3.Range("A2").Value = 45025.99999999999
$ws3.Range("B2").Value = 47
$ws3.Range("C2").Value = -15.5386069732541
$ws3.Range("D2").Value = 112.9575244043992
$ws3.Range("A3").Value = 45046.99999999999
$ws3.Range("B3").Value = 47
$ws3.Range("C3").Value = -19.24704525360733
$ws3.Range("D3").Value = 107.9705995330653
$ws3.Range("A4").Value = 45053.99999999999
$ws3.Range("B4").Value = 47
$ws3.Range("C4").Value = -14.67523667019298
$ws3.Range("D4").Value = 107.36877368816
$ws3.Range("A5").Value = 45060.99999999999
$ws3.Range("B5").Value = 47
$ws3.Range("C5").Value = -16.43559427654041
$ws3.Range("D5").Value = 106.4177911138166
$ws3.Range("A6").Value = 45067.99999999999
$ws3.Range("B6").Value = 47
$ws3.Range("C6").Value = -19.60779893867665
$ws3.Range("D6").Value = 106.3564081860391
$ws3.Range("A7").Value = 45074.99999999999
$ws3.Range("B7").Value = 47
$ws3.Range("C7").Value = -14.63120749911693
$ws3.Range("D7").Value = 113.0537838570039
$ws3.Range("A8").Value = 45081.99999999999
$ws3.Range("B8").Value = 47
$ws3.Range("C8").Value = -11.41995347914994
$ws3.Range("D8").Value = 107.4441376295416
$ws3.Range("A9").Value = 45088.99999999999
$ws3.Range("B9").Value = 47
$ws3.Range("C9").Value = -18.53320107900052
$ws3.Range("D9").Value = 105.4877678039122
$ws3.Range("A10").Value = 45095.99999999999
$ws3.Range("B10").Value = 47
$ws3.Range("C10").Value = -18.41660300060996
$ws3.Range("D10").Value = 110.7691817490527
$ws3.Range("A11").Value = 45102.99999999999
$ws3.Range("B11").Value = 47
$ws3.Range("C11").Value = -13.49853756872608
$ws3.Range("D11").Value = 109.3486093625822
$ws3.Range("A12").Value = 45109.99999999999
$ws3.Range("B12").Value = 47
$ws3.Range("C12").Value = -16.66258835235682
$ws3.Range("D12").Value = 111.6557337841612
$ws3.Range("A13").Value = 45116.99999999999
$ws3.Range("B13").Value = 47
$ws3.Range("C13").Value = -19.59217421950725
$ws3.Range("D13").Value = 108.3832379127191
$ws3.Range("A14").Value = 45123.99999999999
$ws3.Range("B14").Value = 46
$ws3.Range("C14").Value = -19.22711641788708
$ws3.Range("D14").Value = 108.1034325681919
$ws3.Range("A15").Value = 45130.99999999999
$ws3.Range("B15").Value = 46
$ws3.Range("C15").Value = -23.11991000709974
$ws3.Range("D15").Value = 108.7858462412001
$ws3.Range("A16").Value = 45144.99999999999
$ws3.Range("B16").Value = 46
$ws3.Range("C16").Value = -20.69214704694719
$ws3.Range("D16").Value = 114.5794594517367
$ws3.Range("A17").Value = 45158.99999999999
$ws3.Range("B17").Value = 46
$ws3.Range("C17").Value = -17.39519948898724
$ws3.Range("D17").Value = 109.9091195746004
$ws3.Range("A18").Value = 45165.99999999999
$ws3.Range("B18").Value = 46
$ws3.Range("C18").Value = -15.0123033702108
$ws3.Range("D18").Value = 107.9546094042931
$ws3.Range("A19").Value = 45200.99999999999
$ws3.Range("B19").Value = 46
$ws3.Range("C19").Value = -15.4256150086343
$ws3.Range("D19").Value = 109.3186734504024
$ws3.Range("A20").Value = 45214.99999999999
$ws3.Range("B20").Value = 46
$ws3.Range("C20").Value = -19.52124426059091
$ws3.Range("D20").Value = 107.6093537081921
$ws3.Range("A21").Value = 45221.99999999999
$ws3.Range("B21").Value = 46
$ws3.Range("C21").Value = -22.83333386689473
$ws3.Range("D21").Value = 110.1233686501998
$ws3.Range("A22").Value = 45228.99999999999
$ws3.Range("B22").Value = 45
$ws3.Range("C22").Value = -16.40260169531941
$ws3.Range("D22").Value = 111.7374542776398
$ws3.Range("A23").Value = 45263.99999999999
$ws3.Range("B23").Value = 45
$ws3.Range("C23").Value = -22.07543153048564
$ws3.Range("D23").Value = 107.0086827764792
$ws3.Range("A24").Value = 45270.99999999999
$ws3.Range("B24").Value = 45
$ws3.Range("C24").Value = -24.55725495379976
$ws3.Range("D24").Value = 108.4320640317136
$ws3.Range("A25").Value = 45277.99999999999
$ws3.Range("B25").Value = 45
$ws3.Range("C25").Value = -16.48902623034794
$ws3.Range("D25").Value = 106.2699030625866
$ws3.Range("A26").Value = 45298.99999999999
$ws3.Range("B26").Value = 45
$ws3.Range("C26").Value = -23.83406419595261
$ws3.Range("D26").Value = 113.6335991246604
$ws3.Range("A27").Value = 45305.99999999999
$ws3.Range("B27").Value = 45
$ws3.Range("C27").Value = -21.80449768375315
$ws3.Range("D27").Value = 108.760068836617
$ws3.Range("A28").Value = 45312.99999999999
$ws3.Range("B28").Value = 45
$ws3.Range("C28").Value = -19.60473305437578
$ws3.Range("D28").Value = 108.1718475044398
$ws3.Range("A29").Value = 45326.99999999999
$ws3.Range("B29").Value = 45
$ws3.Range("C29").Value = -18.33362235967701
$ws3.Range("D29").Value = 105.1445575118513
$ws3.Range("A30").Value = 45333.99999999999
$ws3.Range("B30").Value = 44
$ws3.Range("C30").Value = -17.16800632145485
$ws3.Range("D30").Value = 102.958468984275
$ws3.Range("A31").Value = 45347.99999999999
$ws3.Range("B31").Value = 44
$ws3.Range("C31").Value = -24.03301593000862
$ws3.Range("D31").Value = 110.6305201681061
$ws3.Range("A32").Value = 45354.99999999999
$ws3.Range("B32").Value = 44
$ws3.Range("C32").Value = -18.05162271451
$ws3.Range("D32").Value = 108.5349926781943
$ws3.Range("A33").Value = 45361.99999999999
$ws3.Range("B33").Value = 44
$ws3.Range("C33").Value = -17.8263348249544
$ws3.Range("D33").Value = 106.6516258686141
$ws3.Range("A34").Value = 45375.99999999999
$ws3.Range("B34").Value = 44
$ws3.Range("C34").Value = -23.38814161794375
$ws3.Range("D34").Value = 110.4996898091871
$ws3.Range("A35").Value = 45382.99999999999
$ws3.Range("B35").Value = 44
$ws3.Range("C35").Value = -22.56784859809259
$ws3.Range("D35").Value = 106.7682429791415
$ws3.Range("A36").Value = 45389.99999999999
$ws3.Range("B36").Value = 44
$ws3.Range("C36").Value = -20.66389732636511
$ws3.Range("D36").Value = 106.0553951221821
$ws3.Range("A37").Value = 45396.99999999999
$ws3.Range("B37").Value = 44
$ws3.Range("C37").Value = -20.9010819005611
$ws3.Range("D37").Value = 107.8971289910267
$ws3.Range("A38").Value = 45515.99999999999
$ws3.Range("B38").Value = 43
$ws3.Range("C38").Value = -19.46072763901211
$ws3.Range("D38").Value = 104.0092137878037
$ws3.Range("A39").Value = 45529.99999999999
$ws3.Range("B39").Value = 43
$ws3.Range("C39").Value = -21.03795212276936
$ws3.Range("D39").Value = 108.2240466498885
$ws3.Range("A40").Value = 45543.99999999999
$ws3.Range("B40").Value = 42
$ws3.Range("C40").Value = -23.97016446564671
$ws3.Range("D40").Value = 107.4039207125151
$ws3.Range("A41").Value = 45550.99999999999
$ws3.Range("B41").Value = 42
$ws3.Range("C41").Value = -24.52990251081908
$ws3.Range("D41").Value = 107.0189256195133
$ws3.Range("A42").Value = 45557.99999999999
$ws3.Range("B42").Value = 42
$ws3.Range("C42").Value = -21.44229238277135
$ws3.Range("D42").Value = 105.8962822628363
$ws3.Range("A43").Value = 45571.99999999999
$ws3.Range("B43").Value = 42
$ws3.Range("C43").Value = -20.71880165877505
$ws3.Range("D43").Value = 103.8445122586859
$ws3.Range("A44").Value = 45578.99999999999
$ws3.Range("B44").Value = 42
$ws3.Range("C44").Value = -22.03349345688114
$ws3.Range("D44").Value = 104.9446191884181
$ws3.Range("A45").Value = 45585.99999999999
$ws3.Range("B45").Value = 42
$ws3.Range("C45").Value = -19.40334565546622
$ws3.Range("D45").Value = 108.1278752843098
$ws3.Range("A46").Value = 45592.99999999999
$ws3.Range("B46").Value = 42
$ws3.Range("C46").Value = -23.48288923051979
$ws3.Range("D46").Value = 106.7178119226458
$ws3.Range("A47").Value = 45599.99999999999
$ws3.Range("B47").Value = 42
$ws3.Range("C47").Value = -18.85305256672553
$ws3.Range("D47").Value = 106.9637981184635
$ws3.Range("A48").Value = 45606.99999999999
$ws3.Range("B48").Value = 42
$ws3.Range("C48").Value = -24.44471503525594
$ws3.Range("D48").Value = 105.6211240349869
$ws3.Range("A49").Value = 45613.99999999999
$ws3.Range("B49").Value = 42
$ws3.Range("C49").Value = -18.72561047330317
$ws3.Range("D49").Value = 106.4696840150491
$ws3.Range("A50").Value = 45620.99999999999
$ws3.Range("B50").Value = 42
$ws3.Range("C50").Value = -23.37880706338751
$ws3.Range("D50").Value = 104.1878207020818
$ws3.Range("A51").Value = 45627.99999999999
$ws3.Range("B51").Value = 42
$ws3.Range("C51").Value = -21.77821846536979
$ws3.Range("D51").Value = 107.7057379168826
$ws3.Range("A52").Value = 45634.99999999999
$ws3.Range("B52").Value = 42
$ws3.Range("C52").Value = -22.02360332243084
$ws3.Range("D52").Value = 106.1436131743291
$ws3.Range("A53").Value = 45641.99999999999
$ws3.Range("B53").Value = 41
$ws3.Range("C53").Value = -20.3342333736896
$ws3.Range("D53").Value = 103.9361693715126
$ws3.Range("A54").Value = 45648.99999999999
$ws3.Range("B54").Value = 41
$ws3.Range("C54").Value = -17.03881079774121
$ws3.Range("D54").Value = 106.7097945664641
$ws3.Range("A55").Value = 45655.99999999999
$ws3.Range("B55").Value = 41
$ws3.Range("C55").Value = -20.06607580847176
$ws3.Range("D55").Value = 102.9111442166208

# copy the date number format from an existing date cell onto the new ds column
$ws1.Range("A2").Copy()
$ws3.Range("A2:A55").PasteSpecial(-4122)

$excel.CutCopyMode = 0
